$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "A" column style (bold/border/centered) used by the ranking
# number cells down into the newly added rows (29-49) before writing values,
# mirroring the style already applied to A2:A28.
$ws.Range("A28").Copy($ws.Range("A29:A49"))

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 'Hakeem Butler'
$ws.Cells.Item(2, 3).Value = 158
$ws.Cells.Item(2, 4).Value = 9.875

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 'J.J. Arcega-Whiteside'
$ws.Cells.Item(3, 3).Value = 135.16
$ws.Cells.Item(3, 4).Value = 8.4475

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 'Darius Slayton'
$ws.Cells.Item(4, 3).Value = 95.8
$ws.Cells.Item(4, 4).Value = 5.9875

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 'Mecole Hardman'
$ws.Cells.Item(5, 3).Value = 95.05
$ws.Cells.Item(5, 4).Value = 5.940625

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 'Jamal Custis'
$ws.Cells.Item(6, 3).Value = 90.83333333333333
$ws.Cells.Item(6, 4).Value = 5.677083333333333

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 'Damarkus Lodge'
$ws.Cells.Item(7, 3).Value = 90.6
$ws.Cells.Item(7, 4).Value = 5.6625

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 'Johnnie Dixon'
$ws.Cells.Item(8, 3).Value = 89.9090909090909
$ws.Cells.Item(8, 4).Value = 5.619318181818182

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 'Deebo Samuel'
$ws.Cells.Item(9, 3).Value = 88.125
$ws.Cells.Item(9, 4).Value = 5.5078125

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 'Felton Davis'
$ws.Cells.Item(10, 3).Value = 85.47619047619048
$ws.Cells.Item(10, 4).Value = 5.342261904761905

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 'Gary Jennings'
$ws.Cells.Item(11, 3).Value = 83.5
$ws.Cells.Item(11, 4).Value = 5.21875

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 'Nyqwan Murray'
$ws.Cells.Item(12, 3).Value = 78.75
$ws.Cells.Item(12, 4).Value = 4.921875

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 'D.K. Metcalf'
$ws.Cells.Item(13, 3).Value = 78.66666666666667
$ws.Cells.Item(13, 4).Value = 4.916666666666667

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 'Riley Ridley'
$ws.Cells.Item(14, 3).Value = 78.16666666666667
$ws.Cells.Item(14, 4).Value = 4.885416666666667

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 'Keesean Johnson'
$ws.Cells.Item(15, 3).Value = 78.16666666666667
$ws.Cells.Item(15, 4).Value = 4.885416666666667

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 'Terry McLaurin'
$ws.Cells.Item(16, 3).Value = 75.05263157894737
$ws.Cells.Item(16, 4).Value = 4.690789473684211

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 'David Sills'
$ws.Cells.Item(17, 3).Value = 72.8
$ws.Cells.Item(17, 4).Value = 4.55

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 'Hunter Renfrow'
$ws.Cells.Item(18, 3).Value = 67.5
$ws.Cells.Item(18, 4).Value = 4.21875

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 'Parris Campbell'
$ws.Cells.Item(19, 3).Value = 66.33333333333333
$ws.Cells.Item(19, 4).Value = 4.145833333333333

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 'Travis Fulgham'
$ws.Cells.Item(20, 3).Value = 49.81818181818182
$ws.Cells.Item(20, 4).Value = 3.113636363636364

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 'Andy Isabella'
$ws.Cells.Item(21, 3).Value = 49.35294117647059
$ws.Cells.Item(21, 4).Value = 3.084558823529412

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 'Miles Boykin'
$ws.Cells.Item(22, 3).Value = 48.4
$ws.Cells.Item(22, 4).Value = 3.025

$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 'Stanley Morgan'
$ws.Cells.Item(23, 3).Value = 47.42857142857143
$ws.Cells.Item(23, 4).Value = 2.964285714285714

$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 'Emmanuel Butler'
$ws.Cells.Item(24, 3).Value = 46
$ws.Cells.Item(24, 4).Value = 2.875

$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = 'Diontae Johnson'
$ws.Cells.Item(25, 3).Value = 40.07142857142857
$ws.Cells.Item(25, 4).Value = 2.504464285714286

$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 'Jalen Hurd'
$ws.Cells.Item(26, 3).Value = 39.52941176470588
$ws.Cells.Item(26, 4).Value = 2.470588235294118

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 'Keelan Doss'
$ws.Cells.Item(27, 3).Value = 34
$ws.Cells.Item(27, 4).Value = 2.125

$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 'Ashton Dulin'
$ws.Cells.Item(28, 3).Value = 33.875
$ws.Cells.Item(28, 4).Value = 2.1171875

$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 'Tyre Brady'
$ws.Cells.Item(29, 3).Value = 32
$ws.Cells.Item(29, 4).Value = 2

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 'Antoine Wesley'
$ws.Cells.Item(30, 3).Value = 28
$ws.Cells.Item(30, 4).Value = 1.75

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 'Kelvin Harmon'
$ws.Cells.Item(31, 3).Value = 28
$ws.Cells.Item(31, 4).Value = 1.75

$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = 'Emanuel Hall'
$ws.Cells.Item(32, 3).Value = 25.3
$ws.Cells.Item(32, 4).Value = 1.58125

$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = 'N''Keal Harry'
$ws.Cells.Item(33, 3).Value = 22.9
$ws.Cells.Item(33, 4).Value = 1.43125

$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = 'Jakobi Meyers'
$ws.Cells.Item(34, 3).Value = 21.33333333333333
$ws.Cells.Item(34, 4).Value = 1.333333333333333

$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = 'A.J. Brown'
$ws.Cells.Item(35, 3).Value = 20.25
$ws.Cells.Item(35, 4).Value = 1.265625

$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = 'Terry Godwin'
$ws.Cells.Item(36, 3).Value = 19.66666666666667
$ws.Cells.Item(36, 4).Value = 1.229166666666667

$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = 'Bisi Johnson'
$ws.Cells.Item(37, 3).Value = 16.8
$ws.Cells.Item(37, 4).Value = 1.05

$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = 'Cody Thompson'
$ws.Cells.Item(38, 3).Value = 13.28571428571429
$ws.Cells.Item(38, 4).Value = 0.8303571428571429

$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = 'Lil''Jordan Humphrey'
$ws.Cells.Item(39, 3).Value = 6
$ws.Cells.Item(39, 4).Value = 0.375

$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = 'Anthony Johnson'
$ws.Cells.Item(40, 3).Value = 4.5
$ws.Cells.Item(40, 4).Value = 0.28125

$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = 'Alex Wesley'
$ws.Cells.Item(41, 3).Value = 3.5
$ws.Cells.Item(41, 4).Value = 0.21875

$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = 'Jaylen Smith'
$ws.Cells.Item(42, 3).Value = 2.5
$ws.Cells.Item(42, 4).Value = 0.15625

$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = 'Jovon Durante'
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(43, 4).Value = 0.0625

$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = 'Ryan Davis'
$ws.Cells.Item(44, 3).Value = 1
$ws.Cells.Item(44, 4).Value = 0.0625

$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).Value = 'Dillon Mitchell'
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(45, 4).Value = 0.0625

$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 'Marquise Brown'
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 0

$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 'Jamarius Way'
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 0

$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 'Greg Dortch'
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 0

$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 'Jazz Ferguson'
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 0

Write-Host "Done updating rows 2-49"